$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1657.5438
$ws.Range("I15").Value = 1657.5438
$ws.Range("K15").Value = 4972.6314
$ws.Range("M15").Value = -4803.6314
$ws.Range("H100").Value = 4057.625
$ws.Range("J100").Value = 6429.25
$ws.Range("L100").Value = 6429.25
$ws.Range("N100").Value = -7511.25
$ws.Range("H116").Value = 6194.5557
$ws.Range("I116").Value = 4707.3335
$ws.Range("J116").Value = 9169
$ws.Range("K116").Value = 4707.3335
$ws.Range("L116").Value = 9169
$ws.Range("M116").Value = -1265.3335
$ws.Range("N116").Value = -16053
$ws.Range("H132").Value = 1564.1666
$ws.Range("I132").Value = 1265.5625
$ws.Range("K132").Value = 3796.6875
$ws.Range("M132").Value = -1266.6875
$ws.Range("H138").Value = 2601.3765
$ws.Range("I138").Value = 1158.6129
$ws.Range("J138").Value = 3322.758
$ws.Range("K138").Value = 3475.8387
$ws.Range("L138").Value = 9968.273999999999
$ws.Range("M138").Value = 1664.1613
$ws.Range("N138").Value = -20248.274
$ws.Range("H141").Value = 4350
$ws.Range("I141").Value = 4350
$ws.Range("K141").Value = 13050
$ws.Range("M141").Value = -7870

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4161.6294
$ws.Range("I2").Value = 359.72726
$ws.Range("J2").Value = 20890
$ws.Range("K2").Value = 359.72726
$ws.Range("L2").Value = 20890
$ws.Range("M2").Value = -246.72726
$ws.Range("N2").Value = -21116
$ws.Range("H32").Value = 2516.1924
$ws.Range("I32").Value = 1746.9459
$ws.Range("K32").Value = 1746.9459
$ws.Range("M32").Value = -1459.9459
$ws.Range("H110").Value = 2370.3635
$ws.Range("I110").Value = 1707.2413
$ws.Range("K110").Value = 1707.2413
$ws.Range("M110").Value = 337.7587000000001
$ws.Range("H116").Value = 4161.6294
$ws.Range("I116").Value = 359.72726
$ws.Range("J116").Value = 20890
$ws.Range("K116").Value = 359.72726
$ws.Range("L116").Value = 20890
$ws.Range("M116").Value = 1934.27274
$ws.Range("N116").Value = -25478
$ws.Range("H122").Value = 3438.348
$ws.Range("I122").Value = 3086
$ws.Range("J122").Value = 3822.7273
$ws.Range("K122").Value = 9258
$ws.Range("L122").Value = 11468.1819
$ws.Range("M122").Value = -6808
$ws.Range("N122").Value = -16368.1819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4161.6294
$ws.Range("I3").Value = 359.72726
$ws.Range("J3").Value = 20890
$ws.Range("K3").Value = 359.72726
$ws.Range("L3").Value = 20890
$ws.Range("M3").Value = -245.72726
$ws.Range("N3").Value = -21118
$ws.Range("H106").Value = 30863.666
$ws.Range("J106").Value = 30863.666
$ws.Range("L106").Value = 30863.666
$ws.Range("N106").Value = -33387.666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 143.9
$ws.Range("I7").Value = 50.363636
$ws.Range("K7").Value = 50.363636
$ws.Range("M7").Value = 62.636364
$ws.Range("H56").Value = 25813.5
$ws.Range("I56").Value = 1255
$ws.Range("J56").Value = 33999.668
$ws.Range("K56").Value = 1255
$ws.Range("L56").Value = 33999.668
$ws.Range("M56").Value = -410
$ws.Range("N56").Value = -35689.668
$ws.Range("H132").Value = 5113.8335
$ws.Range("I132").Value = 4050.5
$ws.Range("J132").Value = 7240.5
$ws.Range("K132").Value = 12151.5
$ws.Range("L132").Value = 21721.5
$ws.Range("M132").Value = -9621.5
$ws.Range("N132").Value = -26781.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 183.6
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""
$ws.Range("H30").Value = 183.6
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = ""
$ws.Range("H113").Value = 1547.8
$ws.Range("J113").Value = 1991.6
$ws.Range("L113").Value = 5974.799999999999
$ws.Range("N113").Value = -10314.8
$ws.Range("H119").Value = 11853.444
$ws.Range("I119").Value = 5324.5
$ws.Range("J119").Value = 13718.857
$ws.Range("K119").Value = 15973.5
$ws.Range("L119").Value = 41156.571
$ws.Range("M119").Value = -11135.5
$ws.Range("N119").Value = -50832.571

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 12000000
$ws.Range("J11").Value = 13000000
$ws.Range("L11").Value = 13000000
$ws.Range("N11").Value = -13000278
$ws.Range("H35").Value = 25015
$ws.Range("I35").Value = 25015
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 25015
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -24717
$ws.Range("N35").Value = ""
$ws.Range("H80").Value = 241645.28
$ws.Range("I80").Value = 387057.38
$ws.Range("J80").Value = 5350.625
$ws.Range("K80").Value = 387057.38
$ws.Range("L80").Value = 5350.625
$ws.Range("M80").Value = -386059.38
$ws.Range("N80").Value = -7346.625
$ws.Range("H83").Value = 241645.28
$ws.Range("I83").Value = 387057.38
$ws.Range("J83").Value = 5350.625
$ws.Range("K83").Value = 1935286.9
$ws.Range("L83").Value = 26753.125
$ws.Range("M83").Value = -1930294.9
$ws.Range("N83").Value = -36737.125
$ws.Range("H122").Value = 8170.931
$ws.Range("I122").Value = 9641.941000000001
$ws.Range("K122").Value = 28925.823
$ws.Range("M122").Value = -26475.823

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14149.5
$ws.Range("I7").Value = 7733.3335
$ws.Range("J7").Value = 17999.2
$ws.Range("K7").Value = 7733.3335
$ws.Range("L7").Value = 17999.2
$ws.Range("M7").Value = -7621.3335
$ws.Range("N7").Value = -18223.2
$ws.Range("H22").Value = 3709.0303
$ws.Range("I22").Value = 1645.8667
$ws.Range("J22").Value = 5428.3335
$ws.Range("K22").Value = 1645.8667
$ws.Range("L22").Value = 5428.3335
$ws.Range("M22").Value = -1350.8667
$ws.Range("N22").Value = -6018.3335
$ws.Range("H27").Value = 3709.0303
$ws.Range("I27").Value = 1645.8667
$ws.Range("J27").Value = 5428.3335
$ws.Range("K27").Value = 1645.8667
$ws.Range("L27").Value = 5428.3335
$ws.Range("M27").Value = -1538.8667
$ws.Range("N27").Value = -5642.3335
$ws.Range("H40").Value = 16312.8
$ws.Range("I40").Value = 19851.334
$ws.Range("J40").Value = 11005
$ws.Range("K40").Value = 19851.334
$ws.Range("L40").Value = 11005
$ws.Range("M40").Value = -19715.334
$ws.Range("N40").Value = -11277
$ws.Range("H61").Value = 3220.76
$ws.Range("I61").Value = 2713
$ws.Range("J61").Value = 5021
$ws.Range("K61").Value = 2713
$ws.Range("L61").Value = 5021
$ws.Range("M61").Value = -2511
$ws.Range("N61").Value = -5425
$ws.Range("H113").Value = 3220.76
$ws.Range("I113").Value = 2713
$ws.Range("J113").Value = 5021
$ws.Range("K113").Value = 2713
$ws.Range("L113").Value = 5021
$ws.Range("M113").Value = -543
$ws.Range("N113").Value = -9361
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""
$ws.Range("H122").Value = 229443.5
$ws.Range("I122").Value = 447396.22
$ws.Range("K122").Value = 1342188.66
$ws.Range("M122").Value = -1339738.66
$ws.Range("H126").Value = 14149.5
$ws.Range("I126").Value = 7733.3335
$ws.Range("J126").Value = 17999.2
$ws.Range("K126").Value = 23200.0005
$ws.Range("L126").Value = 53997.60000000001
$ws.Range("M126").Value = -20730.0005
$ws.Range("N126").Value = -58937.60000000001
$ws.Range("H136").Value = 6043.9375
$ws.Range("I136").Value = 4811.9697
$ws.Range("J136").Value = 8754.267
$ws.Range("K136").Value = 14435.9091
$ws.Range("L136").Value = 26262.801
$ws.Range("M136").Value = -11885.9091
$ws.Range("N136").Value = -31362.801

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = ""
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = ""
